# Added email functionality :-)
#
# - FindAndBookFlight: add a new "closeBrowser" keyword row (Sr No 6) and
#   grow its table to cover it.
# - TestSuite: flip the FindAndBookFlight run flag from N to Y.
$wb = $excel.ActiveWorkbook

# --- FindAndBookFlight sheet: add a new keyword row (closeBrowser) ---
$wsFlight = $wb.Worksheets.Item("FindAndBookFlight")
$wsFlight.Cells.Item(7, 1).Value = 6
$wsFlight.Cells.Item(7, 2).Value = "closeBrowser"

# Grow the sheet's table so the new row is included (Table1423: A1:E6 -> A1:E7)
$loFlight = $wsFlight.ListObjects.Item(1)
$loFlight.Resize($wsFlight.Range("A1:E7"))

$wsFlight.Range("B7").Select()

# --- TestSuite sheet: mark FindAndBookFlight's run flag Y, move selection to A2 ---
$wsTestSuite = $wb.Worksheets.Item("TestSuite")
$wsTestSuite.Activate()
$wsTestSuite.Range("C4").Value = "Y"
$wsTestSuite.Range("A2").Select()
